$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 48619.693414
$ws.Range("D2").Value = 78.198396

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 21367.679467
$ws.Range("D3").Value = 17.183554
$ws.Range("E3").Value = 0

# Row 4 - Residuals
$ws.Range("B4").Value = 206420.324723
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -9.857059
$ws.Range("H5").Value = -18.897534
$ws.Range("I5").Value = -0.816585
$ws.Range("J5").Value = 0.028756

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 7.799861
$ws.Range("H6").Value = -1.695702
$ws.Range("I6").Value = 17.295425
$ws.Range("J6").Value = 0.130823

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 17.65692
$ws.Range("H7").Value = 10.516891
$ws.Range("I7").Value = 24.79695
$ws.Range("J7").Value = 0
